$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.628.51"
$ws.Range("E2").Value = "  -0.48%  "

$ws.Range("D3").Value = "3.679.70"
$ws.Range("E3").Value = "  +2.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.92%  "

$ws.Range("E6").Value = "  +9.64%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "664.24"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.424"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.05%  "

$ws.Range("E9").Value = "  +2.03%  "

$ws.Range("E10").Value = "  +0.01%  "

$ws.Range("D11").Value = "3.677.33"
$ws.Range("E11").Value = "  +2.15%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.88%  "

$ws.Range("E13").Value = "  +0.59%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.67%  "

$ws.Range("D15").Value = "4.363.84"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000269"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.90%  "

$ws.Range("D17").Value = "96.310.08"
$ws.Range("E17").Value = "  -0.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +12.49%  "

$ws.Range("D19").Value = "3.673.53"
$ws.Range("E19").Value = "  +2.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.526"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "528.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.44%  "

$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.21%  "

$ws.Range("E29").Value = "  +12.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.63%  "

$ws.Range("E31").Value = "  +0.95%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +16.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.186"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "32.83"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.38%  "

$ws.Range("E36").Value = "  +0.21%  "

$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "641.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.02%  "

$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.595"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.75"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "44.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +33.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.162"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.968"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.65%  "

$ws.Range("E44").Value = "  +8.44%  "

$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.453"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +21.54%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0462"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.85%  "

$ws.Range("E49").Value = "  -0.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.86%  "
